$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("C1").Value = "total_non_icsap"
$ws.Range("D1").Value = "total_icsap"
$ws.Range("E1").Value = "rate_icsap"
$ws.Range("G1").Value = "non_icsap_cost"
$ws.Range("H1").Value = "icsap_cost"
$ws.Range("I1").Value = "rate_icsap_cost"

# Update data rows 2-12

# Row 2
$ws.Range("B2").Value = 10107830
$ws.Range("C2").Value = 9039770
$ws.Range("D2").Value = 1068060
$ws.Range("E2").Value = 10.57
$ws.Range("F2").Value = 12466941277.75
$ws.Range("G2").Value = 11753619575.6
$ws.Range("H2").Value = 713321702.15
$ws.Range("I2").Value = 0.08

# Row 3
$ws.Range("B3").Value = 10086803
$ws.Range("C3").Value = 9056736
$ws.Range("D3").Value = 1030067
$ws.Range("E3").Value = 10.21
$ws.Range("F3").Value = 12853800102.36
$ws.Range("G3").Value = 12129581538
$ws.Range("H3").Value = 724218564.36

# Row 4
$ws.Range("B4").Value = 10080846
$ws.Range("C4").Value = 9069829
$ws.Range("D4").Value = 1011017
$ws.Range("E4").Value = 10.03
$ws.Range("F4").Value = 13128495829.52
$ws.Range("G4").Value = 12382493363.11
$ws.Range("H4").Value = 746002466.41

# Row 5
$ws.Range("B5").Value = 10263047
$ws.Range("C5").Value = 9254209
$ws.Range("D5").Value = 1008838
$ws.Range("E5").Value = 9.83
$ws.Range("F5").Value = 13663536653.4
$ws.Range("G5").Value = 12889404970.95
$ws.Range("H5").Value = 774131682.45

# Row 6
$ws.Range("B6").Value = 10565062
$ws.Range("C6").Value = 9533098
$ws.Range("D6").Value = 1031964
$ws.Range("E6").Value = 9.77
$ws.Range("F6").Value = 14242190463.88
$ws.Range("G6").Value = 13429586259.96
$ws.Range("H6").Value = 812604203.92
$ws.Range("I6").Value = 0.07

# Row 7
$ws.Range("B7").Value = 10952676
$ws.Range("C7").Value = 9888314
$ws.Range("D7").Value = 1064362
$ws.Range("E7").Value = 9.72
$ws.Range("F7").Value = 14978043954.91
$ws.Range("G7").Value = 14104063926.38
$ws.Range("H7").Value = 873980028.53
$ws.Range("I7").Value = 0.07

# Row 8
$ws.Range("B8").Value = 9382348
$ws.Range("C8").Value = 8563590
$ws.Range("D8").Value = 818758
$ws.Range("E8").Value = 8.73
$ws.Range("F8").Value = 15600721335.6
$ws.Range("G8").Value = 14847115881.68
$ws.Range("H8").Value = 753605453.92
$ws.Range("I8").Value = 0.06

# Row 9
$ws.Range("B9").Value = 10349540
$ws.Range("C9").Value = 9548245
$ws.Range("D9").Value = 801295
$ws.Range("E9").Value = 7.74
$ws.Range("F9").Value = 21598917691.29
$ws.Range("G9").Value = 20831029650.92
$ws.Range("H9").Value = 767888040.37

# Row 10
$ws.Range("B10").Value = 11313938
$ws.Range("C10").Value = 10325202
$ws.Range("D10").Value = 988736
$ws.Range("E10").Value = 8.74
$ws.Range("F10").Value = 18344613234.15
$ws.Range("G10").Value = 17268863031.96
$ws.Range("H10").Value = 1075750202.19
$ws.Range("I10").Value = 0.06

# Row 11
$ws.Range("B11").Value = 12170247
$ws.Range("C11").Value = 11112266
$ws.Range("D11").Value = 1057981
$ws.Range("E11").Value = 8.69
$ws.Range("F11").Value = 19945489853.01
$ws.Range("G11").Value = 18748609846
$ws.Range("H11").Value = 1196880007.01

# Row 12
$ws.Range("B12").Value = 13052583
$ws.Range("C12").Value = 11949413
$ws.Range("D12").Value = 1103170
$ws.Range("E12").Value = 8.45
$ws.Range("F12").Value = 22976380653.58
$ws.Range("G12").Value = 21730243070.21
$ws.Range("H12").Value = 1246137583.37
